$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.783.51"
$ws.Range("E2").Value = "  +4.39%  "
$ws.Range("D3").Value = "3.558.07"
$ws.Range("E3").Value = "  +3.65%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.92%  "
$ws.Range("D7").Value = "3.554.34"
$ws.Range("E7").Value = "  +3.70%  "
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +6.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.55%  "
$ws.Range("E12").Value = "  +3.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("E14").Value = "  +2.98%  "
$ws.Range("D15").Value = "4.139.14"
$ws.Range("E15").Value = "  +4.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "608.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.562.78"
$ws.Range("E18").Value = "  +4.07%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "70.912.48"
$ws.Range("E19").Value = "  +4.59%  "
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -14.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("E29").Value = "  +6.39%  "
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "709.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +17.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.73%  "
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.101"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0478"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.143"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.17%  "
$ws.Range("D43").Value = "3.371.51"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.03%  "
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("E51").Value = "  -0.03%  "
